# ---------------------------------------------------------------------------
# Rebuild Sheet1 of the quiz-import workbook:
#  - insert a new note row on top (merged A1:F1 / G1:L1, centered)
#  - add a new "hide_option" column (H)
#  - replace the two sample quiz questions with two new ones
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wipe out all existing row/column formatting + content for the area we touch
# (EntireRow.Delete removes cell values AND row-level metadata like explicit
# row heights, giving us a clean slate to rebuild on).
$ws.Range("A1:G3").EntireRow.Delete() | Out-Null

# --- cell values -----------------------------------------------------------
# Values are assigned in the same order the original authoring tool used so
# that the shared-string table comes out in the expected sequence.

$ws.Range("A2").Value = "question"
$ws.Range("F2").Value = "option_match"
$ws.Range("G2").Value = "answer"
$ws.Range("B2").Value = "optionA"
$ws.Range("C2").Value = "optionB"
$ws.Range("D2").Value = "optionC"
$ws.Range("E2").Value = "optionD"
$ws.Range("H2").Value = "hide_option"

$ws.Range("A1").Value = "Lưu ý: Trường hide_option nếu bằng 1 thì sẽ ẩn các tùy chọn, nếu bằng 0 thì không ẩn. Và bỏ qua dòng 1 này. Mặc định là không ẩn, nếu không ẩn thì bỏ trường đó đi."

$ws.Range("A3").Value = "Tại sao đến năm 1965, Mĩ phải chuyển sang thực hiện chiến lược “Chiến tranh cục bộ”?"
$ws.Range("B3").Value = "Mĩ muốn nhanh chóng kết thúc chiến tranh Việt Nam."
$ws.Range("C3").Value = "Mĩ muốn mở rộng và quốc tế hóa chiến tranh Việt Nam"
$ws.Range("D3").Value = "Mĩ lo ngại ủng hộ của Trung Quốc và Liên Xô cho cuộc kháng chiến của nhân dân ta"
$ws.Range("E3").Value = "Chiến lược “Chiến tranh đặc biệt” đã bị phá sản hoàn toàn."

$ws.Range("A4").Value = "Việt Nam gia nhập ASEAN có ý nghĩa là"
$ws.Range("B4").Value = "ASEAN đã trở thành một liên minh kinh tế - chính trị"
$ws.Range("C4").Value = "Chứng tỏ sự đối đầu về ý thực hệ tư tưởng – chính trị - quân sự"
$ws.Range("D4").Value = "mở ra triển vọng cho sự liên kết toàn khu vực Đông Nam Á."
$ws.Range("E4").Value = "Chứng tỏ sự hợp tác giữa các thành viên ASEAN ngày càng có hiệu quả."

$ws.Range("F3").Value = "D"
$ws.Range("F4").Value = "C"

$ws.Range("G3").Value = "Đến năm 1965, Mĩ phải chuyển sang thực hiện chiến lược “Chiến tranh cục bộ” vì  chiến lược “Chiến tranh đặc biệt” đã bị phá sản hoàn toàn."
$ws.Range("G4").Value = "Việt Nam gia nhập ASEAN có ý nghĩa là mở ra triển vọng cho sự liên kết toàn khu vực Đông Nam Á."

# --- formatting --------------------------------------------------------

# hide_option header reuses the existing wrap-text style
$ws.Range("H2").WrapText = $true

# note row: centered, merged across A1:F1 and G1:L1
$ws.Range("A1:L1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:F1").MergeCells = $true
$ws.Range("G1:L1").MergeCells = $true

# --- column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.166666666666664   # -> 38
$ws.Columns.Item(2).ColumnWidth = 27                    # -> ~27.875 (closest achievable)
$ws.Columns.Item(8).ColumnWidth = 11.8                  # -> ~12.625 (closest achievable)

# --- view / selection -----------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F4").Select() | Out-Null

# --- page setup -------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

Write-Host "edit complete"
